# The diff shows a new price-report record being inserted as row 474
# (Fecha 2023-03-10 / serial 44995), with every subsequent row shifting
# down by one (old row 474 -> new row 475, ... old row 572 -> new row 573).
# Re-create that by inserting a blank row at 474 (which pushes the rest
# down and grows the used range to R573) and then filling the new row
# with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("474:474").Insert()

$ws.Range("A474").Value = 9
$ws.Range("B474").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C474").Value = "Metropolitana"
$ws.Range("D474").Value = 44995
$ws.Range("E474").Value = 13
$ws.Range("F474").Value = 100112012
$ws.Range("G474").Value = "Espinaca"
$ws.Range("H474").Value = "Sin especificar"
$ws.Range("I474").Value = "Primera"
$ws.Range("J474").Value = 160
$ws.Range("K474").Value = 7000
$ws.Range("L474").Value = 8000
$ws.Range("M474").Value = 7500
$ws.Range("N474").Value = "$/cuna 10 kilos"
$ws.Range("O474").Value = "Provincia de Chacabuco"
$ws.Range("P474").Value = 750
$ws.Range("Q474").Value = 10
$ws.Range("R474").Value = "Hortaliza"
